$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then a hashtable of Column -> new value.
# D (Price) and E (Volume 1h) columns hold numeric-looking text that must
# stay literal text (e.g. "2.130", "39.00"), so those cells are forced to
# the Text number format before the value is written, otherwise Excel would
# auto-convert them to numbers/percentages and drop significant trailing zeros.
$updates = @(
    @{ Row = 2; Cells = @{ "D"="313.78"; "E"="2.78%" } }
    @{ Row = 3; Cells = @{ "D"="35.19"; "E"="-2.07%" } }
    @{ Row = 4; Cells = @{ "D"="5.115"; "E"="0.80%" } }
    @{ Row = 5; Cells = @{ "D"="0.08153"; "E"="2.82%" } }
    @{ Row = 6; Cells = @{ "D"="2.130"; "E"="0.31%" } }
    @{ Row = 7; Cells = @{ "B"="KuCoinToken"; "C"="https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"; "D"="7.959"; "E"="-0.03%" } }
    @{ Row = 8; Cells = @{ "B"="MXToken"; "C"="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; "D"="0.9298"; "E"="0.68%" } }
    @{ Row = 9; Cells = @{ "B"="LiechtensteinCryptoassetsExchange"; "C"="https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; "D"="0.1029"; "E"="5.76%" } }
    @{ Row = 10; Cells = @{ "B"="WazirX"; "C"="https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; "D"="0.1945"; "E"="4.92%" } }
    @{ Row = 11; Cells = @{ "B"="MandalaExchangeToken"; "C"="https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; "D"="0.08981"; "E"="4.52%" } }
    @{ Row = 12; Cells = @{ "B"="BitrueCoin"; "C"="https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; "D"="0.03716"; "E"="4.19%" } }
    @{ Row = 13; Cells = @{ "B"="BitMartToken"; "C"="https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; "D"="0.09916"; "E"="-0.26%" } }
    @{ Row = 14; Cells = @{ "B"="BitForexToken"; "C"="https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; "D"="0.001444"; "E"="-0.03%" } }
    @{ Row = 15; Cells = @{ "B"="TigerCash"; "C"="https://coinranking.com/coin/6hIn06L2+tigercash-tch"; "D"="0.005726"; "E"="0.06%" } }
    @{ Row = 16; Cells = @{ "B"="LEO"; "C"="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; "D"="3.469"; "E"="0.17%" } }
    @{ Row = 17; Cells = @{ "B"="GateToken"; "C"="https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; "D"="4.154"; "E"="0.52%" } }
    @{ Row = 18; Cells = @{ "D"="2.972"; "E"="8.08%" } }
    @{ Row = 19; Cells = @{ "D"="0.3413"; "E"="0.77%" } }
    @{ Row = 20; Cells = @{ "E"="-1.27%" } }
    @{ Row = 21; Cells = @{ "D"="5.104"; "E"="-1.38%" } }
    @{ Row = 22; Cells = @{ "E"="-0.07%" } }
    @{ Row = 23; Cells = @{ "D"="0.04574"; "E"="0.43%" } }
    @{ Row = 24; Cells = @{ "D"="0.001249"; "E"="0.93%" } }
    @{ Row = 25; Cells = @{ "D"="0.004699"; "E"="-3.81%" } }
    @{ Row = 26; Cells = @{ "D"="0.0001254"; "E"="-3.82%" } }
    @{ Row = 27; Cells = @{ "D"="0.0004508"; "E"="-5.31%" } }
    @{ Row = 39; Cells = @{ "D"="0.01955"; "E"="5.24%" } }
    @{ Row = 40; Cells = @{ "D"="0.04876"; "E"="2.81%" } }
    @{ Row = 41; Cells = @{ "D"="0.007499"; "E"="-5.32%" } }
    @{ Row = 42; Cells = @{ "D"="0.1390"; "E"="-0.61%" } }
    @{ Row = 43; Cells = @{ "D"="0.007882"; "E"="1.84%" } }
    @{ Row = 44; Cells = @{ "D"="0.002101"; "E"="-4.27%" } }
    @{ Row = 45; Cells = @{ "D"="0.01178"; "E"="4.40%" } }
    @{ Row = 46; Cells = @{ "D"="0.00006756"; "E"="7.61%" } }
    @{ Row = 47; Cells = @{ "D"="0.00000000752"; "E"="0.05%" } }
    @{ Row = 48; Cells = @{ "D"="39.00"; "E"="-21.94%" } }
    @{ Row = 49; Cells = @{ "D"="0.001703"; "E"="-15.05%" } }
    @{ Row = 50; Cells = @{ "D"="0.00002106"; "E"="0.05%" } }
    @{ Row = 51; Cells = @{ "D"="0.0002006"; "E"="0.05%" } }
)

$textFormatColumns = @("D", "E")

foreach ($update in $updates) {
    $row = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $cellRef = "$col$row"
        $cell = $ws.Range($cellRef)
        if ($textFormatColumns -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $update.Cells[$col]
    }
}
